$d = $word.ActiveDocument

# 1) Contact info line: merge the split "github.com/danielmartincraig" run
#    (and surrounding proofErr spell-check tags) with the rest of the line
#    into a single run.
$d.Content.Find.Execute(
    "github.com/danielmartincraig " + [char]0x2022 + " linkedin.com/danielcraig23",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "github.com/danielmartincraig " + [char]0x2022 + " linkedin.com/danielcraig23",
    2) | Out-Null

# 2) Insert a new "OBJECTIVE:" heading paragraph right after the contact
#    info line (i.e. right before the "EDUCATION:" heading).
$edu = $d.Paragraphs(4)
$edu.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs(4)

$objectiveFrag = '<w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">OBJECTIVE: </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Eager to drive back-end solutions at </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>3M</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> on a full-time basis</w:t></w:r></w:p>'

$pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $objectiveFrag + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newPara.Range.InsertXML($pkg) | Out-Null

# 3) "Web Engineering I and II" bullet: merge the split "Engineering" run
#    (and gramStart/gramEnd proofErr tags) into a single run.
$d.Content.Find.Execute(
    [char]0x2022 + "    Web Engineering I and II",
    $true, $false, $false, $false, $false, $true, 1, $false,
    [char]0x2022 + "    Web Engineering I and II",
    2) | Out-Null

# 4) Remove the stray "_GoBack" bookmark that used to sit at the very end
#    of the document (after "Fluent in Spanish").
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

Write-Host "done"
